$wb = $excel.ActiveWorkbook

# Both the "展览" (Exhibitions) sheet and the "全部类型" (All Types) sheet
# contain the same two data rows; update the "想去人数" (F) column values
# for row 2 (96 -> 98) and row 3 (15 -> 17) on each.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 98
    $ws.Range("F3").Value = 17
}
